# Tutorial9.pptx edit: update the Assignment 3 deadline time on the
# "Assignment Overview" slide from "00:00:00 am" to "18:00:00 pm".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldText = "DDL: 00:00:00 am, Apr 17"
$newText = "DDL: 18:00:00 pm, Apr 17"

$full = $tr.Text
$idx = $full.IndexOf($oldText)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}
